$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target text for D2 (Abstract) and E2 (Authors) were base64-encoded (UTF-8)
# to guarantee exact fidelity of whitespace/newlines/unicode characters
# when reverting these two cells to their original shared-string content.
$d2Base64 = "UHVycG9zZQppZD0iUGFyMSI+VG8gYXNzZXNzIHRoZSBtYW5hZ2VtZW50IGFuZCBzYWZldHkgb2YgZXBpZHVyYWwgb3IgZ2VuZXJhbCBhbmVzdGhlc2lhIGZvciBDZXNhcmVhbiBkZWxpdmVyeSBpbiBwYXJ0dXJpZW50cyB3aXRoIGNvcm9uYXZpcnVzIGRpc2Vhc2UgKENPVklELTE5KSBhbmQgdGhlaXIgbmV3Ym9ybnMsIGFuZCB0byBldmFsdWF0ZSB0aGUgc3RhbmRhcmRpemVkIHByb2NlZHVyZXMgZm9yIHByb3RlY3RpbmcgbWVkaWNhbCBzdGFmZi4KCgpNZXRob2RzCmlkPSJQYXIyIj5XZSByZXRyb3NwZWN0aXZlbHkgcmV2aWV3ZWQgdGhlIGNhc2VzIG9mIHBhcnR1cmllbnRzIGRpYWdub3NlZCB3aXRoIHNldmVyZSBhY3V0ZSByZXNwaXJhdG9yeSBzeW5kcm9tZSBjb3JvbmF2aXJ1cyAoU0FSUy1Db1YtMikgaW5mZWN0aW9uIGRpc2Vhc2UgKENPVklELTE5KS4KCiBUaGVpciBlcGlkZW1pb2xvZ2ljIGhpc3RvcnksIGNoZXN0IGNvbXB1dGVkIHRvbW9ncmFwaHkgc2NhbnMsIGxhYm9yYXRvcnkgbWVhc3VyZW1lbnRzLCBhbmQgU0FSUy1Db1YtMiBudWNsZWljIGFjaWQgcG9zaXRpdml0eSB3ZXJlIGV2YWx1YXRlZC4KCiBXZSBhbHNvIHJlY29yZGVkIHRoZSBwYXRpZW50c+KAmSBkZW1vZ3JhcGhpYyBhbmQgY2xpbmljYWzCoGNoYXJhY3RlcmlzdGljcywgYW5lc3RoZXNpYSBhbmQgc3VyZ2VyeS1yZWxhdGVkIGRhdGEsIG1hdGVybmFsIGFuZCBuZW9uYXRhbCBjb21wbGljYXRpb25zLCBhcyB3ZWxsIGFzIHRoZSBoZWFsdGggc3RhdHVzIG9mIHRoZSBpbnZvbHZlZCBtZWRpY2FsIHN0YWZmLgoKClJlc3VsdHMKaWQ9IlBhcjMiPlRoZSBjbGluaWNhbCBjaGFyYWN0ZXJpc3RpY3Mgb2YgMTcgcHJlZ25hbnQgd29tZW4gaW5mZWN0ZWQgd2l0aCBTQVJTLUNvVi0yIHdlcmUgc2ltaWxhciB0byB0aG9zZSBwcmV2aW91c2x5IHJlcG9ydGVkIGluIG5vbi1wcmVnbmFudCBhZHVsdCBwYXRpZW50cy4KCiBBbGwgb2YgdGhlIDE3IHBhdGllbnRzIHVuZGVyd2VudCBDZXNhcmVhbiBkZWxpdmVyeSB3aXRoIGFuZXN0aGVzaWEgcGVyZm9ybWVkIGFjY29yZGluZyB0byBzdGFuZGFyZGl6ZWQgYW5lc3RoZXNpYS9zdXJnZXJ5IHByb2NlZHVyZXMuCgogRm91cnRlZW4gb2YgdGhlIHBhdGllbnRzIHVuZGVyd2VudCBjb250aW51b3VzIGVwaWR1cmFsIGFuZXN0aGVzaWEgd2l0aCAxMiBleHBlcmllbmNpbmcgc2lnbmlmaWNhbnQgaW50cmFvcGVyYXRpdmUgaHlwb3RlbnNpb24uCgogVGhyZWUgcGF0aWVudHMgcmVjZWl2ZWQgZ2VuZXJhbCBhbmVzdGhlc2lhIHdpdGggdHJhY2hlYWwgaW50dWJhdGlvbiBiZWNhdXNlIGVtZXJnZW5jeSBzdXJnZXJ5IHdhcyBuZWVkZWQuCgogVGhyZWUgb2YgdGhlIHBhcnR1cmllbnRzIGFyZSBzdGlsbCByZWNvdmVyaW5nIGZyb20gdGhlaXIgQ2VzYXJlYW4gZGVsaXZlcnkgYW5kIGFyZSByZWNlaXZpbmcgaW4taG9zcGl0YWwgdHJlYXRtZW50IGZvciBDT1ZJRC0xOS4gVGhyZWUgbmVvbmF0ZXMgd2VyZSBib3JuIHByZW1hdHVyZWx5LgoKIFRoZXJlIHdlcmUgbm8gZGVhdGhzIG9yIHNlcmlvdXMgbmVvbmF0YWwgYXNwaHl4aWEgZXZlbnRzLgoKIEFsbCBuZW9uYXRhbCBTQVJTLUNvVi0yIG51Y2xlaWMgYWNpZCB0ZXN0cyB3ZXJlIG5lZ2F0aXZlLgoKIE5vIG1lZGljYWwgc3RhZmYgd2VyZSBpbmZlY3RlZCB0aHJvdWdob3V0IHRoZSBwYXRpZW50IGNhcmUgcGVyaW9kLgoKCkNvbmNsdXNpb25zCmlkPSJQYXI0Ij5Cb3RoIGVwaWR1cmFsIGFuZCBnZW5lcmFsIGFuZXN0aGVzaWEgd2VyZSBzYWZlbHkgdXNlZCBmb3IgQ2VzYXJlYW4gZGVsaXZlcnkgaW4gdGhlIHBhcnR1cmllbnRzIHdpdGggQ09WSUQtMTkuIE5ldmVydGhlbGVzcywgdGhlIGluY2lkZW5jZSBvZiBoeXBvdGVuc2lvbiBkdXJpbmcgZXBpZHVyYWwgYW5lc3RoZXNpYSBhcHBlYXJlZCBleGNlc3NpdmUuCgogUHJvcGVyIHBhdGllbnQgdHJhbnNmZXIsIG1lZGljYWwgc3RhZmYgYWNjZXNzIHByb2NlZHVyZXMsIGFuZCBlZmZlY3RpdmUgYmlvc2FmZXR5IHByZWNhdXRpb25zIGFyZSBpbXBvcnRhbnQgdG8gcHJvdGVjdCBtZWRpY2FsIHN0YWZmIGZyb20gQ09WSUQtMTkuCg=="
$e2Base64 = "W1JvbmclQ2hlbiVOVUxMJTEsIFl1YW4lWmhhbmclTlVMTCUxLCBMZWklSHVhbmclTlVMTCUxLCBCaS1oZW5nJUNoZW5nJU5VTEwlMSwgWmhvbmcteXVhbiVYaWElTlVMTCUxLCBRaW5nLXRhbyVNZW5nJW1lbmdxaW5ndGFvMjAxOEAxMjYuY29tJTFd"

$d2Bytes = [System.Convert]::FromBase64String($d2Base64)
$e2Bytes = [System.Convert]::FromBase64String($e2Base64)

$d2Text = [System.Text.Encoding]::UTF8.GetString($d2Bytes)
$e2Text = [System.Text.Encoding]::UTF8.GetString($e2Bytes)

$ws.Range("D2").Value2 = $d2Text
$ws.Range("E2").Value2 = $e2Text
